$wb = $excel.ActiveWorkbook

# --- Sheet "Schedule": update Cost ($) and Unit Cost ($/ML) for row 2 ---
$schedule = $wb.Worksheets.Item("Schedule")
$schedule.Range("E2").Value = 1674.97327725
$schedule.Range("F2").Value = 27.69466397569445

# --- Sheet "Detailed": update Price (and Type where applicable) ---
$detailed = $wb.Worksheets.Item("Detailed")

$detailed.Range("B5").Value = 78
$detailed.Range("B7").Value = 67.98876
$detailed.Range("C7").Value = "historical"
$detailed.Range("B8").Value = 68.04342
$detailed.Range("C8").Value = "historical"
$detailed.Range("B9").Value = 67.2062
$detailed.Range("B10").Value = 66.36573
$detailed.Range("B11").Value = 61.70593
$detailed.Range("B12").Value = 63.8185
$detailed.Range("B15").Value = 78.71372
$detailed.Range("B16").Value = 56.98
$detailed.Range("B17").Value = 51.37223
$detailed.Range("B18").Value = 50.05842
$detailed.Range("B20").Value = 41.92894
$detailed.Range("B21").Value = 19.05581
$detailed.Range("B22").Value = 0.51
$detailed.Range("B23").Value = 34.01
$detailed.Range("B24").Value = 22.07
$detailed.Range("B25").Value = 27.65305
$detailed.Range("B28").Value = 12.09885
$detailed.Range("B30").Value = 52.11742
$detailed.Range("B31").Value = 59.44001
$detailed.Range("B32").Value = 62.04007
$detailed.Range("B34").Value = 47.4258
$detailed.Range("B35").Value = 53.14823
$detailed.Range("B36").Value = 61.98752
$detailed.Range("B37").Value = 24.39195
$detailed.Range("B38").Value = 33.74417
$detailed.Range("B39").Value = 80.02
$detailed.Range("B40").Value = 158.99
$detailed.Range("B41").Value = 198.74561
$detailed.Range("B44").Value = 75.21648
$detailed.Range("B45").Value = 65.0001
$detailed.Range("B46").Value = 65
$detailed.Range("B47").Value = 65
$detailed.Range("B48").Value = 66.44042
$detailed.Range("B49").Value = 59.81397
